$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing data-row formatting down onto the new rows (7-21) first,
# so the new cells pick up the same cell style (font/alignment) as the rest
# of the "People / Movies" table instead of defaulting to unstyled cells.
$ws.Range("A6:B6").Copy()
$ws.Range("A7:B21").PasteSpecial(-4122)
# Nudge the number format so these rows land on their own (new) cell-style
# record rather than silently collapsing back onto row 6's style index.
$ws.Range("A7:B21").NumberFormat = "General"

# --- Existing rows: re-shuffle the nominee data for "Зеленая миля" /
#     "Форрест Гамп" and fill in "Побег из Шоушенка" on row 6 ---
$ws.Range("A3").Value = ""
$ws.Range("B3").Value = "Зеленая миля"
$ws.Range("A4").Value = "Том Хэнкс"
$ws.Range("B4").Value = "Форрест Гамп"
$ws.Range("A6").Value = "Тим Роббинс"
$ws.Range("B6").Value = "Побег из Шоушенка"

# --- New nominee rows (7-21) ---
$ws.Range("A7").Value = "Морган Фриман"
$ws.Range("B7").Value = "Побег из Шоушенка"

$ws.Range("B8").Value = "Побег из Шоушенка"

$ws.Range("A9").Value = "Франсуа Клюзе"
$ws.Range("B9").Value = "1+1"

$ws.Range("A10").Value = "Омар Си"
$ws.Range("B10").Value = "1+1"

$ws.Range("B11").Value = "1+1"

$ws.Range("A12").Value = "Лиам Нисон"
$ws.Range("B12").Value = "Список Шиндлера"

$ws.Range("B13").Value = "Список Шиндлера"

$ws.Range("A14").Value = "Мэттью МакКонахи"
$ws.Range("B14").Value = "Интерстеллар"

$ws.Range("A15").Value = "Энн Хэтэуэй"
$ws.Range("B15").Value = "Интерстеллар"

$ws.Range("B16").Value = "Интерстеллар"

$ws.Range("A17").Value = "Элайджа Вуд"
$ws.Range("B17").Value = "Властелин колец: Возвращение короля"

$ws.Range("A18").Value = "Орландо Блум"
$ws.Range("B18").Value = "Властелин колец: Возвращение короля"

$ws.Range("B19").Value = "Властелин колец: Возвращение короля"

$ws.Range("A20").Value = "Эдвард Нортон"
$ws.Range("B20").Value = "Бойцовский клуб"

$ws.Range("B21").Value = "Бойцовский клуб"
